$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.566.11"
$ws.Range("E2").Value = "  -2.55%  "

$ws.Range("D3").Value = "3.529.19"
$ws.Range("E3").Value = "  -3.53%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.90"
$ws.Range("E5").Value = "  -4.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.06"
$ws.Range("E6").Value = "  -3.24%  "

$ws.Range("D7").Value = "3.526.30"
$ws.Range("E7").Value = "  -3.48%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.485"
$ws.Range("E9").Value = "  -2.03%  "

$ws.Range("E10").Value = "  -2.12%  "

$ws.Range("E11").Value = "  -1.29%  "

$ws.Range("E12").Value = "  -1.14%  "

$ws.Range("E13").Value = "  -2.35%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.11"
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.124.68"
$ws.Range("E15").Value = "  -3.67%  "

$ws.Range("D16").Value = "3.522.33"
$ws.Range("E16").Value = "  -4.43%  "

$ws.Range("D17").Value = "67.638.22"
$ws.Range("E17").Value = "  -2.45%  "

$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.47"
$ws.Range("E20").Value = "  -2.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "452.12"
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.41"
$ws.Range("E22").Value = "  -2.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.639"
$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("E24").Value = "  -1.79%  "

$ws.Range("D25").Value = "3.675.82"
$ws.Range("E25").Value = "  -3.46%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("E27").Value = "  -3.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.49"
$ws.Range("E28").Value = "  -1.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.34"
$ws.Range("E29").Value = "  -5.28%  "

$ws.Range("E30").Value = "  -0.90%  "

$ws.Range("E31").Value = "  +1.12%  "

$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.90"
$ws.Range("E33").Value = "  -1.95%  "

$ws.Range("E34").Value = "  -3.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.20"
$ws.Range("E35").Value = "  -2.75%  "

$ws.Range("E36").Value = "  -2.28%  "

$ws.Range("D37").Value = "3.529.16"
$ws.Range("E37").Value = "  -3.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.95"
$ws.Range("E38").Value = "  -3.55%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.56"
$ws.Range("E41").Value = "  -0.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.59"
$ws.Range("E42").Value = "  -4.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0874"
$ws.Range("E43").Value = "  -1.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.08"
$ws.Range("E44").Value = "  -3.71%  "

$ws.Range("E45").Value = "  -3.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.33"
$ws.Range("E46").Value = "  +10.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.78"
$ws.Range("E47").Value = "  -1.54%  "

$ws.Range("E48").Value = "  -3.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.64"
$ws.Range("E49").Value = "  -1.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.21"
$ws.Range("E50").Value = "  -2.53%  "

$ws.Range("E51").Value = "  -2.26%  "
